$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the "_old" / "_new" header labels to "_FV2210" / "_FV2304" ---
$headers = @(
    "Segmentname_FV2210","Segmentgruppe_FV2210","Segment_FV2210","Datenelement_FV2210",
    "Segment ID_FV2210","Code_FV2210","Qualifier_FV2210","Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210","Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304","Segmentgruppe_FV2304","Segment_FV2304","Datenelement_FV2304",
    "Segment ID_FV2304","Code_FV2304","Qualifier_FV2304","Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304","Bedingung_FV2304"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2) Turn the data range into an Excel Table ("Table1") with autofilter ---
$tableRange = $ws.Range("A1:U79")
$list = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$list.Name = "Table1"

# --- 3) Freeze the header row (pane split below row 1) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$null = $excel.ActiveWindow.FreezePanes
